$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A (rows 2-54): remove 3 duplicate entries, shift remaining values up ---
$ws.Cells.Item(2, 1).Value = 1800
$ws.Cells.Item(3, 1).Value = 1096
$ws.Cells.Item(4, 1).Value = 1095
$ws.Cells.Item(5, 1).Value = 1009
$ws.Cells.Item(6, 1).Value = 1000
$ws.Cells.Item(7, 1).Value = 1008
$ws.Cells.Item(8, 1).Value = 1020
$ws.Cells.Item(9, 1).Value = 1030
$ws.Cells.Item(10, 1).Value = 1315
$ws.Cells.Item(11, 1).Value = 1321
$ws.Cells.Item(12, 1).Value = 1322
$ws.Cells.Item(13, 1).Value = 1329
$ws.Cells.Item(14, 1).Value = 1337
$ws.Cells.Item(15, 1).Value = 1317
$ws.Cells.Item(16, 1).Value = 1318
$ws.Cells.Item(17, 1).Value = 1326
$ws.Cells.Item(18, 1).Value = 1323
$ws.Cells.Item(19, 1).Value = 1353
$ws.Cells.Item(20, 1).Value = 1324
$ws.Cells.Item(21, 1).Value = 1356
$ws.Cells.Item(22, 1).Value = 1357
$ws.Cells.Item(23, 1).Value = 1338
$ws.Cells.Item(24, 1).Value = 1422
$ws.Cells.Item(25, 1).Value = 1400
$ws.Cells.Item(26, 1).Value = 1434
$ws.Cells.Item(27, 1).Value = 1373
$ws.Cells.Item(28, 1).Value = 1372
$ws.Cells.Item(29, 1).Value = 1312
$ws.Cells.Item(30, 1).Value = 1302
$ws.Cells.Item(31, 1).Value = 1429
$ws.Cells.Item(32, 1).Value = 1442
$ws.Cells.Item(33, 1).Value = 1441
$ws.Cells.Item(34, 1).Value = 1436
$ws.Cells.Item(35, 1).Value = 1405
$ws.Cells.Item(36, 1).Value = 1432
$ws.Cells.Item(37, 1).Value = 1433
$ws.Cells.Item(38, 1).Value = 1426
$ws.Cells.Item(39, 1).Value = 1375
$ws.Cells.Item(40, 1).Value = 1092
$ws.Cells.Item(41, 1).Value = 1094
$ws.Cells.Item(42, 1).Value = 1091
$ws.Cells.Item(43, 1).Value = 1098
$ws.Cells.Item(44, 1).Value = 1097
$ws.Cells.Item(45, 1).Value = 1350
$ws.Cells.Item(46, 1).Value = 1260
$ws.Cells.Item(47, 1).Value = 1196
$ws.Cells.Item(48, 1).Value = 1180
$ws.Cells.Item(49, 1).Value = 1165
$ws.Cells.Item(50, 1).Value = 1110
$ws.Cells.Item(51, 1).Value = 1162
$ws.Cells.Item(52, 1).Value = 1163
$ws.Cells.Item(53, 1).Value = 1027
$ws.Cells.Item(54, 1).Value = 1028

# --- Clear now-unused tail rows of column A (was 57 rows, now 54) ---
$ws.Range("A55:A57").ClearContents()

# --- Update column B (rows 2-8): remove 1 duplicate entry, shift remaining values up ---
$ws.Cells.Item(2, 2).Value = 1122
$ws.Cells.Item(3, 2).Value = 1037
$ws.Cells.Item(4, 2).Value = 1042
$ws.Cells.Item(5, 2).Value = 1040
$ws.Cells.Item(6, 2).Value = 1043
$ws.Cells.Item(7, 2).Value = 1044
$ws.Cells.Item(8, 2).Value = 1038
$ws.Range("B9").ClearContents()

# --- Update the active cell / selection ---
[void]$ws.Range("D20").Select()

# --- Highlight duplicate values (Conditional Formatting), matching the
#     standard "Light Red Fill with Dark Red Text" built-in style.      ---
$dupFontColor = 393372      # BGR for RGB(0x9C,0x00,0x06)
$dupFillColor = 13551615    # BGR for RGB(0xFF,0xC7,0xCE)

$dupAreas = @("A10:C13", "C2:C9", "B2:B8", "A2:A9", "A14:A54", "B14:C57")
foreach ($addr in $dupAreas) {
    $rng = $ws.Range($addr)
    $fc = $rng.FormatConditions.AddUniqueValues()
    $fc.DupeUnique = 1
    $fc.Font.Color = $dupFontColor
    $fc.Interior.Color = $dupFillColor
}
